$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "29.539.87"
Set-TextValue "E2" "  +0.95%  "
Set-TextValue "D3" "1.983.47"
Set-TextValue "D4" "1.005"
Set-TextValue "E4" "  +0.47%  "
Set-TextValue "D5" "327.83"
Set-TextValue "E5" "  +0.56%  "
Set-TextValue "D6" "1.004"
Set-TextValue "E6" "  +0.45%  "
Set-TextValue "D7" "0.4665"
Set-TextValue "E7" "  +0.40%  "
Set-TextValue "D8" "0.3920"
Set-TextValue "E8" "  +0.23%  "
Set-TextValue "D9" "46.18"
Set-TextValue "E9" "  -1.26%  "
Set-TextValue "D10" "0.07950"
Set-TextValue "E10" "  +0.91%  "
Set-TextValue "D11" "0.9943"
Set-TextValue "E11" "  +0.57%  "
Set-TextValue "D12" "22.86"
Set-TextValue "E12" "  +3.79%  "
Set-TextValue "D13" "1.980.20"
Set-TextValue "E13" "  +2.62%  "
Set-TextValue "D14" "7.196"
Set-TextValue "E14" "  +1.61%  "
Set-TextValue "D15" "5.836"
Set-TextValue "E15" "  +1.53%  "
Set-TextValue "D16" "0.07101"
Set-TextValue "E16" "  +1.62%  "
Set-TextValue "D17" "87.71"
Set-TextValue "E17" "  -0.75%  "
Set-TextValue "D18" "1.007"
Set-TextValue "E18" "  +0.54%  "
Set-TextValue "D19" "0.000009966"
Set-TextValue "E19" "  -0.45%  "
Set-TextValue "D20" "17.33"
Set-TextValue "E20" "  +1.40%  "
Set-TextValue "E21" "  +0.36%  "
Set-TextValue "D22" "29.556.12"
Set-TextValue "E22" "  +0.97%  "
Set-TextValue "D23" "5.574"
Set-TextValue "E23" "  +4.90%  "
Set-TextValue "E24" "  +1.05%  "
Set-TextValue "D25" "2.221.51"
Set-TextValue "E25" "  +3.29%  "
Set-TextValue "D26" "2.109"
Set-TextValue "E26" "  +0.79%  "
Set-TextValue "D27" "158.86"
Set-TextValue "E27" "  +1.60%  "
Set-TextValue "D28" "19.62"
Set-TextValue "E28" "  +0.84%  "
Set-TextValue "D29" "5.858"
Set-TextValue "E29" "  -2.28%  "
Set-TextValue "D30" "119.64"
Set-TextValue "E30" "  +0.86%  "
Set-TextValue "D31" "1.902"
Set-TextValue "E31" "  -0.85%  "
Set-TextValue "D32" "0.09429"
Set-TextValue "E32" "  +0.67%  "
Set-TextValue "D33" "0.8947"
Set-TextValue "E33" "  -1.19%  "
Set-TextValue "D34" "5.238"
Set-TextValue "E34" "  -0.92%  "
Set-TextValue "D35" "1.327"
Set-TextValue "E35" "  +0.03%  "
Set-TextValue "D36" "3.199"
Set-TextValue "E36" "  -0.64%  "
Set-TextValue "D37" "0.05815"
Set-TextValue "E37" "  +0.47%  "
Set-TextValue "D38" "1.177"
Set-TextValue "E38" "  -0.38%  "
Set-TextValue "D39" "0.02098"
Set-TextValue "E39" "  +0.38%  "
Set-TextValue "E40" "  +0.92%  "
Set-TextValue "D41" "0.5742"
Set-TextValue "E41" "  +0.43%  "
Set-TextValue "D42" "0.1808"
Set-TextValue "E42" "  +1.18%  "
Set-TextValue "D43" "0.000003088"
Set-TextValue "E43" "  +34.73%  "
Set-TextValue "D44" "9.696"
Set-TextValue "E44" "  -0.93%  "
Set-TextValue "D45" "2.805"
Set-TextValue "E45" "  +7.86%  "
Set-TextValue "D46" "11.87"
Set-TextValue "E46" "  -0.68%  "
Set-TextValue "D47" "0.5370"
Set-TextValue "E47" "  +0.41%  "
Set-TextValue "D48" "2.176"
Set-TextValue "E48" "  -1.15%  "
Set-TextValue "D49" "0.06947"
Set-TextValue "E49" "  -1.39%  "
Set-TextValue "E50" "  +0.86%  "
Set-TextValue "D51" "1.828"
Set-TextValue "E51" "  -1.66%  "
